$wb = $excel.ActiveWorkbook

# ---- Sheet: "Heap Report from Test" ----
$ws = $wb.Worksheets.Item("Heap Report from Test")
$ws.Activate()

# Row 3: add Description / Observations for the "Full test of PushP and PushGP" run
$ws.Cells.Item(3, 9).Value  = "Full test of PushP and PushGP."
$ws.Cells.Item(3, 10).Value = "Memory usage seems random between 83M and 403M"

# Row 4: new data row (Laptop / Release / Factory_Class / 130095 memory test / notes)

# Copy number formatting from row 3 onto row 4 so the date & numeric columns
# keep the same display format as the rest of the table.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("F3:H3").Copy() | Out-Null
$ws.Range("F4:H4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Cells.Item(4, 1).Value  = 43411.40625            # 11/7/2018 9:45 AM
$ws.Cells.Item(4, 2).Value  = "Laptop"
$ws.Cells.Item(4, 3).Value  = "Release"
$ws.Cells.Item(4, 4).Value  = "Factory_Class"
$ws.Cells.Item(4, 6).Value  = 130095
$ws.Cells.Item(4, 7).Value  = 130095
$ws.Cells.Item(4, 8).Value  = 130095
$ws.Cells.Item(4, 9).Value  = "Limited test to PushP.  Also, updated pack() in Literal.h to use factory."
$ws.Cells.Item(4, 10).Value = "Memory usage consistant"

# Update the view: selection on I4, scrolled so column B is left-most visible column
$ws.Range("I4").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
